$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(70, 8).Value = 2298.6667  # H70: 2198.6667 -> 2298.6667
$ws.Cells.Item(70, 9).Value = 2298.6667  # I70: 2198.6667 -> 2298.6667
$ws.Cells.Item(70, 11).Value = 6896.000100000001  # K70: 6596.000100000001 -> 6896.000100000001
$ws.Cells.Item(70, 13).Value = -6626.000100000001  # M70: -6326.000100000001 -> -6626.000100000001

$ws.Cells.Item(73, 8).Value = 2298.6667  # H73: 2198.6667 -> 2298.6667
$ws.Cells.Item(73, 9).Value = 2298.6667  # I73: 2198.6667 -> 2298.6667
$ws.Cells.Item(73, 11).Value = 6896.000100000001  # K73: 6596.000100000001 -> 6896.000100000001
$ws.Cells.Item(73, 13).Value = -5960.000100000001  # M73: -5660.000100000001 -> -5960.000100000001

$ws.Cells.Item(86, 8).Value = 3048.25  # H86: 3000 -> 3048.25
$ws.Cells.Item(86, 9).Value = 2996.3333  # I86: 3000 -> 2996.3333
$ws.Cells.Item(86, 10).Value = 3204  # J86: 0 -> 3204
$ws.Cells.Item(86, 11).Value = 2996.3333  # K86: 3000 -> 2996.3333
$ws.Cells.Item(86, 12).Value = 3204  # L86: 0 -> 3204
$ws.Cells.Item(86, 13).Value = -1873.3333  # M86: -1877 -> -1873.3333
$ws.Cells.Item(86, 14).Value = -5450  # N86: None -> -5450

$ws.Cells.Item(89, 8).Value = 3048.25  # H89: 3000 -> 3048.25
$ws.Cells.Item(89, 9).Value = 2996.3333  # I89: 3000 -> 2996.3333
$ws.Cells.Item(89, 10).Value = 3204  # J89: 0 -> 3204
$ws.Cells.Item(89, 11).Value = 14981.6665  # K89: 15000 -> 14981.6665
$ws.Cells.Item(89, 12).Value = 16020  # L89: 0 -> 16020
$ws.Cells.Item(89, 13).Value = -9365.666499999999  # M89: -9384 -> -9365.666499999999
$ws.Cells.Item(89, 14).Value = -27252  # N89: None -> -27252

$ws.Cells.Item(107, 8).Value = 180.66667  # H107: 187.5 -> 180.66667
$ws.Cells.Item(107, 9).Value = 168.2  # I107: 175.22223 -> 168.2
$ws.Cells.Item(107, 10).Value = 205.6  # J107: 209.6 -> 205.6
$ws.Cells.Item(107, 11).Value = 168.2  # K107: 175.22223 -> 168.2
$ws.Cells.Item(107, 12).Value = 205.6  # L107: 209.6 -> 205.6
$ws.Cells.Item(107, 13).Value = 1751.8  # M107: 1744.77777 -> 1751.8
$ws.Cells.Item(107, 14).Value = -4045.6  # N107: -4049.6 -> -4045.6

$ws.Cells.Item(125, 8).Value = 313.75  # H125: 315.2 -> 313.75
$ws.Cells.Item(125, 9).Value = 300  # I125: 266 -> 300
$ws.Cells.Item(125, 10).Value = 318.33334  # J125: 348 -> 318.33334
$ws.Cells.Item(125, 11).Value = 2700  # K125: 2394 -> 2700
$ws.Cells.Item(125, 12).Value = 2865.00006  # L125: 3132 -> 2865.00006
$ws.Cells.Item(125, 13).Value = -240  # M125: 66 -> -240
$ws.Cells.Item(125, 14).Value = -7785.00006  # N125: -8052 -> -7785.00006

$ws.Cells.Item(132, 8).Value = 1202  # H132: 1177 -> 1202
$ws.Cells.Item(132, 9).Value = 660.8570999999999  # I132: 700.375 -> 660.8570999999999
$ws.Cells.Item(132, 11).Value = 1982.5713  # K132: 2101.125 -> 1982.5713
$ws.Cells.Item(132, 13).Value = 547.4287000000002  # M132: 428.875 -> 547.4287000000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 670.1  # H5: 617.36365 -> 670.1
$ws.Cells.Item(5, 9).Value = 812.5  # I5: 732.2222 -> 812.5
$ws.Cells.Item(5, 11).Value = 812.5  # K5: 732.2222 -> 812.5
$ws.Cells.Item(5, 13).Value = -700.5  # M5: -620.2222 -> -700.5

$ws.Cells.Item(32, 8).Value = 2838.52  # H32: 2950.7083 -> 2838.52
$ws.Cells.Item(32, 9).Value = 2838.52  # I32: 2950.7083 -> 2838.52
$ws.Cells.Item(32, 11).Value = 2838.52  # K32: 2950.7083 -> 2838.52
$ws.Cells.Item(32, 13).Value = -2551.52  # M32: -2663.7083 -> -2551.52

$ws.Cells.Item(35, 8).Value = 2634.25  # H35: 2307.4 -> 2634.25
$ws.Cells.Item(35, 9).Value = 1845.6666  # I35: 1634.25 -> 1845.6666
$ws.Cells.Item(35, 11).Value = 1845.6666  # K35: 1634.25 -> 1845.6666
$ws.Cells.Item(35, 13).Value = -1439.6666  # M35: -1228.25 -> -1439.6666

$ws.Cells.Item(37, 8).Value = 9994.5  # H37: 13996 -> 9994.5
$ws.Cells.Item(37, 10).Value = 0  # J37: 21999 -> 0
$ws.Cells.Item(37, 12).Value = 0  # L37: 21999 -> 0
$ws.Cells.Item(37, 14).ClearContents()  # N37: remove (was -22545)

$ws.Cells.Item(63, 8).Value = 7072.727  # H63: 8233.75 -> 7072.727
$ws.Cells.Item(63, 9).Value = 5499.25  # I63: 5624.25 -> 5499.25
$ws.Cells.Item(63, 10).Value = 11268.667  # J63: 13452.75 -> 11268.667
$ws.Cells.Item(63, 11).Value = 5499.25  # K63: 5624.25 -> 5499.25
$ws.Cells.Item(63, 12).Value = 11268.667  # L63: 13452.75 -> 11268.667
$ws.Cells.Item(63, 13).Value = -4813.25  # M63: -4938.25 -> -4813.25
$ws.Cells.Item(63, 14).Value = -12640.667  # N63: -14824.75 -> -12640.667

$ws.Cells.Item(66, 8).Value = 7072.727  # H66: 8233.75 -> 7072.727
$ws.Cells.Item(66, 9).Value = 5499.25  # I66: 5624.25 -> 5499.25
$ws.Cells.Item(66, 10).Value = 11268.667  # J66: 13452.75 -> 11268.667
$ws.Cells.Item(66, 11).Value = 27496.25  # K66: 28121.25 -> 27496.25
$ws.Cells.Item(66, 12).Value = 56343.335  # L66: 67263.75 -> 56343.335
$ws.Cells.Item(66, 13).Value = -24064.25  # M66: -24689.25 -> -24064.25
$ws.Cells.Item(66, 14).Value = -63207.335  # N66: -74127.75 -> -63207.335

$ws.Cells.Item(74, 8).Value = 1268.25  # H74: 1258.1428 -> 1268.25
$ws.Cells.Item(74, 9).Value = 1268.25  # I74: 1258.1428 -> 1268.25
$ws.Cells.Item(74, 11).Value = 1268.25  # K74: 1258.1428 -> 1268.25
$ws.Cells.Item(74, 13).Value = -394.25  # M74: -384.1428000000001 -> -394.25

$ws.Cells.Item(77, 8).Value = 1268.25  # H77: 1258.1428 -> 1268.25
$ws.Cells.Item(77, 9).Value = 1268.25  # I77: 1258.1428 -> 1268.25
$ws.Cells.Item(77, 11).Value = 6341.25  # K77: 6290.714 -> 6341.25
$ws.Cells.Item(77, 13).Value = -1973.25  # M77: -1922.714 -> -1973.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 670.1  # H4: 617.36365 -> 670.1
$ws.Cells.Item(4, 9).Value = 812.5  # I4: 732.2222 -> 812.5
$ws.Cells.Item(4, 11).Value = 812.5  # K4: 732.2222 -> 812.5
$ws.Cells.Item(4, 13).Value = -697.5  # M4: -617.2222 -> -697.5

$ws.Cells.Item(12, 8).Value = 80  # H12: 490 -> 80
$ws.Cells.Item(12, 9).Value = 0  # I12: 900 -> 0
$ws.Cells.Item(12, 11).Value = 0  # K12: 900 -> 0
$ws.Cells.Item(12, 13).ClearContents()  # M12: remove (was -732)

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(133, 8).Value = 50000  # H133: 0 -> 50000
$ws.Cells.Item(133, 10).Value = 50000  # J133: 0 -> 50000
$ws.Cells.Item(133, 12).Value = 50000  # L133: 0 -> 50000
$ws.Cells.Item(133, 14).Value = -60120  # N133: None -> -60120

$ws.Cells.Item(16, 8).Value = 2248.5  # H16: 2624.25 -> 2248.5
$ws.Cells.Item(16, 9).Value = 2165.6667  # I16: 2499 -> 2165.6667
$ws.Cells.Item(16, 10).Value = 2497  # J16: 3000 -> 2497
$ws.Cells.Item(16, 11).Value = 2165.6667  # K16: 2499 -> 2165.6667
$ws.Cells.Item(16, 12).Value = 2497  # L16: 3000 -> 2497
$ws.Cells.Item(16, 13).Value = -1878.6667  # M16: -2212 -> -1878.6667
$ws.Cells.Item(16, 14).Value = -3071  # N16: -3574 -> -3071

$ws.Cells.Item(31, 8).Value = 2144.3428  # H31: 2219.3235 -> 2144.3428
$ws.Cells.Item(31, 9).Value = 1193.9584  # I31: 1263.4783 -> 1193.9584
$ws.Cells.Item(31, 11).Value = 1193.9584  # K31: 1263.4783 -> 1193.9584
$ws.Cells.Item(31, 13).Value = -898.9584  # M31: -968.4783 -> -898.9584

$ws.Cells.Item(32, 8).Value = 0  # H32: 5270.3335 -> 0
$ws.Cells.Item(32, 9).Value = 0  # I32: 2900 -> 0
$ws.Cells.Item(32, 10).Value = 0  # J32: 10011 -> 0
$ws.Cells.Item(32, 11).Value = 0  # K32: 2900 -> 0
$ws.Cells.Item(32, 12).Value = 0  # L32: 10011 -> 0
$ws.Cells.Item(32, 13).ClearContents()  # M32: remove (was -2584)
$ws.Cells.Item(32, 14).ClearContents()  # N32: remove (was -10643)

$ws.Cells.Item(34, 8).Value = 2144.3428  # H34: 2219.3235 -> 2144.3428
$ws.Cells.Item(34, 9).Value = 1193.9584  # I34: 1263.4783 -> 1193.9584
$ws.Cells.Item(34, 11).Value = 1193.9584  # K34: 1263.4783 -> 1193.9584
$ws.Cells.Item(34, 13).Value = -991.9584  # M34: -1061.4783 -> -991.9584

$ws.Cells.Item(113, 8).Value = 2248.5  # H113: 2624.25 -> 2248.5
$ws.Cells.Item(113, 9).Value = 2165.6667  # I113: 2499 -> 2165.6667
$ws.Cells.Item(113, 10).Value = 2497  # J113: 3000 -> 2497
$ws.Cells.Item(113, 11).Value = 2165.6667  # K113: 2499 -> 2165.6667
$ws.Cells.Item(113, 12).Value = 2497  # L113: 3000 -> 2497
$ws.Cells.Item(113, 13).Value = 4.333299999999781  # M113: -329 -> 4.333299999999781
$ws.Cells.Item(113, 14).Value = -6837  # N113: -7340 -> -6837

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 1499.5  # H5: 1500 -> 1499.5
$ws.Cells.Item(5, 10).Value = 1499.5  # J5: 1500 -> 1499.5
$ws.Cells.Item(5, 12).Value = 4498.5  # L5: 4500 -> 4498.5
$ws.Cells.Item(5, 14).Value = -4722.5  # N5: -4724 -> -4722.5

$ws.Cells.Item(23, 8).Value = 225  # H23: 248 -> 225
$ws.Cells.Item(23, 9).Value = 100  # I23: 0 -> 100
$ws.Cells.Item(23, 10).Value = 266.66666  # J23: 248 -> 266.66666
$ws.Cells.Item(23, 11).Value = 300  # K23: 0 -> 300
$ws.Cells.Item(23, 12).Value = 799.9999799999999  # L23: 744 -> 799.9999799999999
$ws.Cells.Item(23, 13).Value = -65  # M23: None -> -65
$ws.Cells.Item(23, 14).Value = -1269.99998  # N23: -1214 -> -1269.99998

$ws.Cells.Item(135, 8).Value = 1499.5  # H135: 1500 -> 1499.5
$ws.Cells.Item(135, 10).Value = 1499.5  # J135: 1500 -> 1499.5
$ws.Cells.Item(135, 12).Value = 13495.5  # L135: 13500 -> 13495.5
$ws.Cells.Item(135, 14).Value = -18565.5  # N135: -18570 -> -18565.5

$ws.Cells.Item(137, 8).Value = 2874.75  # H137: 4166.6665 -> 2874.75
$ws.Cells.Item(137, 9).Value = 2499.5  # I137: 2500 -> 2499.5
$ws.Cells.Item(137, 10).Value = 3250  # J137: 5000 -> 3250
$ws.Cells.Item(137, 11).Value = 7498.5  # K137: 7500 -> 7498.5
$ws.Cells.Item(137, 12).Value = 9750  # L137: 15000 -> 9750
$ws.Cells.Item(137, 13).Value = -2398.5  # M137: -2400 -> -2398.5
$ws.Cells.Item(137, 14).Value = -19950  # N137: -25200 -> -19950

$ws.Cells.Item(140, 8).Value = 5861.4  # H140: 6744.25 -> 5861.4
$ws.Cells.Item(140, 9).Value = 5861.4  # I140: 6744.25 -> 5861.4
$ws.Cells.Item(140, 11).Value = 17584.2  # K140: 20232.75 -> 17584.2
$ws.Cells.Item(140, 13).Value = -12404.2  # M140: -15052.75 -> -12404.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 1875  # H102: 2081.5 -> 1875
$ws.Cells.Item(102, 9).Value = 1428.5714  # I102: 1497.8 -> 1428.5714
$ws.Cells.Item(102, 11).Value = 1428.5714  # K102: 1497.8 -> 1428.5714
$ws.Cells.Item(102, 13).Value = 193.4286  # M102: 124.2 -> 193.4286

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(32, 8).Value = 24999  # H32: 13 -> 24999
$ws.Cells.Item(32, 9).Value = 24999  # I32: 13 -> 24999
$ws.Cells.Item(32, 11).Value = 24999  # K32: 13 -> 24999
$ws.Cells.Item(32, 13).Value = -24682  # M32: 304 -> -24682

$ws.Cells.Item(55, 8).Value = 843.3333  # H55: 943.75 -> 843.3333
$ws.Cells.Item(55, 9).Value = 766.6667  # I55: 912 -> 766.6667
$ws.Cells.Item(55, 11).Value = 766.6667  # K55: 912 -> 766.6667
$ws.Cells.Item(55, 13).Value = -593.6667  # M55: -739 -> -593.6667

$ws.Cells.Item(132, 8).Value = 2609.9412  # H132: 2809.8823 -> 2609.9412
$ws.Cells.Item(132, 9).Value = 2226.5  # I132: 2469.2856 -> 2226.5
$ws.Cells.Item(132, 11).Value = 6679.5  # K132: 7407.8568 -> 6679.5
$ws.Cells.Item(132, 13).Value = -4149.5  # M132: -4877.8568 -> -4149.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 777.3333  # H126: 819.4 -> 777.3333
$ws.Cells.Item(126, 9).Value = 818.8  # I126: 819.4 -> 818.8
$ws.Cells.Item(126, 10).Value = 570  # J126: 0 -> 570
$ws.Cells.Item(126, 11).Value = 2456.4  # K126: 2458.2 -> 2456.4
$ws.Cells.Item(126, 12).Value = 1710  # L126: 0 -> 1710
$ws.Cells.Item(126, 13).Value = 13.60000000000036  # M126: 11.80000000000018 -> 13.60000000000036
$ws.Cells.Item(126, 14).Value = -6650  # N126: None -> -6650

$ws.Cells.Item(136, 8).Value = 933  # H136: 873.375 -> 933
$ws.Cells.Item(136, 9).Value = 933  # I136: 872.8333 -> 933
$ws.Cells.Item(136, 10).Value = 0  # J136: 875 -> 0
$ws.Cells.Item(136, 11).Value = 2799  # K136: 2618.4999 -> 2799
$ws.Cells.Item(136, 12).Value = 0  # L136: 2625 -> 0
$ws.Cells.Item(136, 13).Value = -249  # M136: -68.4998999999998 -> -249
$ws.Cells.Item(136, 14).ClearContents()  # N136: remove (was -7725)
